$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.986.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -5.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.220.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.34%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.68%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.23%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.591'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.29%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.562'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.95%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.58%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.84%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0832'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.12%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.54%  '

# Row 14
$ws.Range("E14").Value = '  -2.51%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.862'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -11.93%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.557.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.32%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.63%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.220.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.45%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.887.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.29%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0963'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.23%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -11.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.49%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -10.39%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '235.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.88%  '

# Row 26
$ws.Range("E26").Value = '  -6.92%  '

# Row 27
$ws.Range("E27").Value = '  +0.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.77%  '

# Row 29
$ws.Range("E29").Value = '  -4.75%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.33'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -12.46%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0893'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.35%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.62%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.20%  '

# Row 35
$ws.Range("E35").Value = '  -7.59%  '

# Row 36
$ws.Range("E36").Value = '  +9.44%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.70%  '

# Row 38
$ws.Range("E38").Value = '  -6.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.22%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.107'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.82%  '

# Row 41
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.59%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0325'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.97%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.917.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.11%  '

# Row 44
$ws.Range("E44").Value = '  +0.14%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.38%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.75%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.208'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.13%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.74%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '77.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '60.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -12.93%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.63%  '
